$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI TPM-derived ligand-receptor values for L1cam-Itga5 pairs.
# Each entry maps a cell address to its new numeric value, as produced by
# rerunning the analysis with the new TPM expression matrix.
$updates = @{
    "G2" = 6.603177
    "H2" = 19.809531
    "I2" = 0.5135477412645301
    "J2" = 0.5135477412645302
    "M2" = 35.04689966666667
    "N2" = 105.140699
    "O2" = 0.3824629895491901
    "P2" = 0.3824629895491901
    "Q2" = 231.420881800241
    "R2" = 2082.787936202169
    "S2" = 0.1964130044002662
    "T2" = 0.1964130044002662
    "G3" = 6.603177
    "H3" = 19.809531
    "I3" = 0.5135477412645301
    "J3" = 0.5135477412645302
    "O3" = 0.3264402385872224
    "P3" = 0.3264402385872223
    "Q3" = 197.522609855613
    "R3" = 1777.703488700517
    "S3" = 0.1676426471843223
    "T3" = 0.1676426471843223
    "G4" = 6.603177
    "H4" = 19.809531
    "I4" = 0.5135477412645301
    "J4" = 0.5135477412645302
    "M4" = 8.911727666666666
    "N4" = 26.735183
    "O4" = 0.09725271102035077
    "P4" = 0.09725271102035075
    "Q4" = 58.84571515879699
    "R4" = 529.611436429173
    "S4" = 0.04994391007635321
    "T4" = 0.04994391007635321
    "G5" = 6.603177
    "H5" = 19.809531
    "I5" = 0.5135477412645301
    "J5" = 0.5135477412645302
    "M5" = 17.76285166666667
    "N5" = 53.288555
    "O5" = 0.1938440608432367
    "P5" = 0.1938440608432367
    "Q5" = 117.291253579745
    "R5" = 1055.621282217705
    "S5" = 0.09954817960358836
    "T5" = 0.09954817960358837
    "I6" = 0.02944398858046029
    "J6" = 0.0294439885804603
    "M6" = 35.04689966666667
    "N6" = 105.140699
    "O6" = 0.3824629895491901
    "P6" = 0.3824629895491901
    "Q6" = 13.26839406250344
    "R6" = 119.415546562531
    "S6" = 0.01126123589673506
    "T6" = 0.01126123589673506
    "I7" = 0.02944398858046029
    "J7" = 0.0294439885804603
    "O7" = 0.3264402385872224
    "P7" = 0.3264402385872223
    "S7" = 0.009611702657164909
    "T7" = 0.009611702657164909
    "I8" = 0.02944398858046029
    "J8" = 0.0294439885804603
    "M8" = 8.911727666666666
    "N8" = 26.735183
    "O8" = 0.09725271102035077
    "P8" = 0.09725271102035075
    "Q8" = 3.373888006747444
    "R8" = 30.364992060727
    "S8" = 0.002863507712702013
    "T8" = 0.002863507712702013
    "I9" = 0.02944398858046029
    "J9" = 0.0294439885804603
    "M9" = 17.76285166666667
    "N9" = 53.288555
    "O9" = 0.1938440608432367
    "P9" = 0.1938440608432367
    "Q9" = 6.724832091532777
    "R9" = 60.52348882379501
    "S9" = 0.005707542313858313
    "T9" = 0.005707542313858312
    "G10" = 3.441487333333333
    "H10" = 10.324462
    "I10" = 0.2676541983690312
    "J10" = 0.2676541983690313
    "M10" = 35.04689966666667
    "N10" = 105.140699
    "O10" = 0.3824629895491901
    "P10" = 0.3824629895491901
    "Q10" = 120.6134612754376
    "R10" = 1085.521151478938
    "S10" = 0.1023678248736116
    "T10" = 0.1023678248736117
    "G11" = 3.441487333333333
    "H11" = 10.324462
    "I11" = 0.2676541983690312
    "J11" = 0.2676541983690313
    "O11" = 0.3264402385872224
    "P11" = 0.3264402385872223
    "Q11" = 102.9461363620927
    "R11" = 926.515227258834
    "S11" = 0.08737310037445829
    "T11" = 0.08737310037445829
    "G12" = 3.441487333333333
    "H12" = 10.324462
    "I12" = 0.2676541983690312
    "J12" = 0.2676541983690313
    "M12" = 8.911727666666666
    "N12" = 26.735183
    "O12" = 0.09725271102035077
    "P12" = 0.09725271102035075
    "Q12" = 30.66959788294955
    "R12" = 276.026380946546
    "S12" = 0.02603009640736703
    "T12" = 0.02603009640736704
    "G13" = 3.441487333333333
    "H13" = 10.324462
    "I13" = 0.2676541983690312
    "J13" = 0.2676541983690313
    "M13" = 17.76285166666667
    "N13" = 53.288555
    "O13" = 0.1938440608432367
    "P13" = 0.1938440608432367
    "Q13" = 61.13062901471222
    "R13" = 550.17566113241
    "S13" = 0.05188317671359424
    "T13" = 0.05188317671359424
    "G14" = 2.434707333333333
    "H14" = 7.304122
    "I14" = 0.1893540717859783
    "J14" = 0.1893540717859783
    "M14" = 35.04689966666667
    "N14" = 105.140699
    "O14" = 0.3824629895491901
    "P14" = 0.3824629895491901
    "Q14" = 85.32894362903089
    "R14" = 767.9604926612782
    "S14" = 0.07242092437857721
    "T14" = 0.07242092437857721
    "G15" = 2.434707333333333
    "H15" = 7.304122
    "I15" = 0.1893540717859783
    "J15" = 0.1893540717859783
    "O15" = 0.3264402385872224
    "P15" = 0.3264402385872223
    "Q15" = 72.83005539827266
    "R15" = 655.470498584454
    "S15" = 0.06181278837127679
    "T15" = 0.06181278837127679
    "G16" = 2.434707333333333
    "H16" = 7.304122
    "I16" = 0.1893540717859783
    "J16" = 0.1893540717859783
    "M16" = 8.911727666666666
    "N16" = 26.735183
    "O16" = 0.09725271102035077
    "P16" = 0.09725271102035075
    "Q16" = 21.69744870270289
    "R16" = 195.277038324326
    "S16" = 0.0184151968239285
    "T16" = 0.0184151968239285
    "G17" = 2.434707333333333
    "H17" = 7.304122
    "I17" = 0.1893540717859783
    "J17" = 0.1893540717859783
    "M17" = 17.76285166666667
    "N17" = 53.288555
    "O17" = 0.1938440608432367
    "P17" = 0.1938440608432367
    "Q17" = 43.24734521374555
    "R17" = 389.2261069237101
    "S17" = 0.0367051622121958
    "T17" = 0.0367051622121958
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
